# This workbook contains weekly price data for "Betarraga" (beet) at the
# Vega Central Mapocho de Santiago market. A new week of data (2 rows: one
# "Primera" quality row and one "Segunda" quality row) was appended to the
# historical series. Because the data appears to be sorted with the newest
# week inserted right after the existing row 818 block (pushing the
# previously-last weeks down), we insert two new rows at 819:820 - this
# shifts all the old rows 819-931 down to 821-933 intact (including every
# column), and then we just need to populate the two freshly inserted rows
# with the new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows, pushing old rows 819-931 down to 821-933.
$ws.Rows("819:820").Insert()

# New row 819: Betarraga, "Primera" quality entry for the new date.
$ws.Range("A819").Value = 9
$ws.Range("B819").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C819").Value = "Metropolitana"
$ws.Range("D819").Value2 = 45127
$ws.Range("E819").Value = 13
$ws.Range("F819").Value = 100114014
$ws.Range("G819").Value = "Betarraga"
$ws.Range("H819").Value = "Sin especificar"
$ws.Range("I819").Value = "Primera"
$ws.Range("J819").Value = 9700
$ws.Range("K819").Value = 100
$ws.Range("L819").Value = 110
$ws.Range("M819").Value = 105
$ws.Range("N819").Value = "$/unidad"
$ws.Range("O819").Value = "Región Metropolitana"
$ws.Range("P819").Value = 105
$ws.Range("Q819").Value = 1
$ws.Range("R819").Value = "Hortaliza"

# New row 820: Betarraga, "Segunda" quality entry for the same new date.
$ws.Range("A820").Value = 9
$ws.Range("B820").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C820").Value = "Metropolitana"
$ws.Range("D820").Value2 = 45127
$ws.Range("E820").Value = 13
$ws.Range("F820").Value = 100114014
$ws.Range("G820").Value = "Betarraga"
$ws.Range("H820").Value = "Sin especificar"
$ws.Range("I820").Value = "Segunda"
$ws.Range("J820").Value = 7000
$ws.Range("K820").Value = 80
$ws.Range("L820").Value = 80
$ws.Range("M820").Value = 80
$ws.Range("N820").Value = "$/unidad"
$ws.Range("O820").Value = "Región Metropolitana"
$ws.Range("P820").Value = 80
$ws.Range("Q820").Value = 1
$ws.Range("R820").Value = "Hortaliza"
